# Fill in the "খাতা/পত্রের সংখ্যা" (quantity) inputs for a few billing rows
# so the dependent per-row amount formulas (and the grand-total SUM) recalc.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G16").Value = 27   # row 16: I16 = IF(H16=0,0,G16*K16/H16)
$ws.Range("G20").Value = 40   # row 20: I20 = IF(G20=0,0,G20*K20)
$ws.Range("G26").Value = 1    # row 26: I26 = K26*G26

# I32 (grand total, =SUM(I9:I31)) picks up the change automatically on recalc.
